$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = -21.602
$ws.Range("D11").Value = -7.262
$ws.Range("A12").Value = -21.544
$ws.Range("A15").Value = -22.062
$ws.Range("D23").Value = -8.363
$ws.Range("A27").Value = -21.686
$ws.Range("A28").Value = -21.89
$ws.Range("D28").Value = -8.172999999999998
$ws.Range("A31").Value = -21.849
$ws.Range("A32").Value = -21.582
$ws.Range("D32").Value = -7.306999999999999
$ws.Range("D34").Value = -7.933
$ws.Range("A36").Value = -20.673
$ws.Range("D36").Value = -7.681999999999999
$ws.Range("D37").Value = -8.279
$ws.Range("A38").Value = -19.992
$ws.Range("D42").Value = -8.563000000000001
$ws.Range("A46").Value = -21.873
$ws.Range("D49").Value = -8.276
$ws.Range("A54").Value = -22.002
$ws.Range("D54").Value = -7.877000000000001
$ws.Range("A55").Value = -22.184
$ws.Range("A56").Value = -21.933
$ws.Range("A67").Value = -21.577
$ws.Range("A69").Value = -21.503
$ws.Range("A72").Value = -21.689
$ws.Range("A73").Value = -20.078
$ws.Range("D78").Value = -8.382
$ws.Range("D80").Value = -8.270999999999999
$ws.Range("A83").Value = -21.987
$ws.Range("A86").Value = -22.135
$ws.Range("A91").Value = -20.887
$ws.Range("A93").Value = -21.452
$ws.Range("D97").Value = -7.674000000000001
$ws.Range("A99").Value = -22.131
$ws.Range("D99").Value = -8.350999999999999
$ws.Range("D100").Value = -8.324000000000002
$ws.Range("D101").Value = -8.019
$ws.Range("A104").Value = -21.385
$ws.Range("A105").Value = -20.345
